$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.3747116666666667
$ws.Cells.Item(2, 8).Value = 1.124135
$ws.Cells.Item(2, 9).Value = 0.3914669751594584
$ws.Cells.Item(2, 10).Value = 0.3914669751594584
$ws.Cells.Item(2, 13).Value = 7.757543333333333
$ws.Cells.Item(2, 14).Value = 23.27263
$ws.Cells.Item(2, 15).Value = 0.4040769763164727
$ws.Cells.Item(2, 16).Value = 0.4040769763164727
$ws.Cells.Item(2, 17).Value = 2.906841991672223
$ws.Cells.Item(2, 18).Value = 26.16157792505
$ws.Cells.Item(2, 19).Value = 0.1581827916501897
$ws.Cells.Item(2, 20).Value = 0.1581827916501897
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.3747116666666667
$ws.Cells.Item(3, 8).Value = 1.124135
$ws.Cells.Item(3, 9).Value = 0.3914669751594584
$ws.Cells.Item(3, 10).Value = 0.3914669751594584
$ws.Cells.Item(3, 15).Value = 0.01627055103446774
$ws.Cells.Item(3, 16).Value = 0.01627055103446774
$ws.Cells.Item(3, 17).Value = 0.1170468097583334
$ws.Cells.Item(3, 18).Value = 1.053421287825
$ws.Cells.Item(3, 19).Value = 0.006369383397640684
$ws.Cells.Item(3, 20).Value = 0.006369383397640684
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.3747116666666667
$ws.Cells.Item(4, 8).Value = 1.124135
$ws.Cells.Item(4, 9).Value = 0.3914669751594584
$ws.Cells.Item(4, 10).Value = 0.3914669751594584
$ws.Cells.Item(4, 13).Value = 11.12827366666667
$ws.Cells.Item(4, 14).Value = 33.384821
$ws.Cells.Item(4, 15).Value = 0.5796524726490594
$ws.Cells.Item(4, 16).Value = 0.5796524726490595
$ws.Cells.Item(4, 17).Value = 4.169893972759446
$ws.Cells.Item(4, 18).Value = 37.529045754835
$ws.Cells.Item(4, 19).Value = 0.226914800111628
$ws.Cells.Item(4, 20).Value = 0.226914800111628
$ws.Cells.Item(5, 9).Value = 0.4195022558883632
$ws.Cells.Item(5, 10).Value = 0.4195022558883631
$ws.Cells.Item(5, 13).Value = 7.757543333333333
$ws.Cells.Item(5, 14).Value = 23.27263
$ws.Cells.Item(5, 15).Value = 0.4040769763164727
$ws.Cells.Item(5, 16).Value = 0.4040769763164727
$ws.Cells.Item(5, 17).Value = 3.115018252870001
$ws.Cells.Item(5, 18).Value = 28.03516427583
$ws.Cells.Item(5, 19).Value = 0.169511203117309
$ws.Cells.Item(5, 20).Value = 0.169511203117309
$ws.Cells.Item(6, 9).Value = 0.4195022558883632
$ws.Cells.Item(6, 10).Value = 0.4195022558883631
$ws.Cells.Item(6, 15).Value = 0.01627055103446774
$ws.Cells.Item(6, 16).Value = 0.01627055103446774
$ws.Cells.Item(6, 19).Value = 0.006825532863505958
$ws.Cells.Item(6, 20).Value = 0.006825532863505957
$ws.Cells.Item(7, 9).Value = 0.4195022558883632
$ws.Cells.Item(7, 10).Value = 0.4195022558883631
$ws.Cells.Item(7, 13).Value = 11.12827366666667
$ws.Cells.Item(7, 14).Value = 33.384821
$ws.Cells.Item(7, 15).Value = 0.5796524726490594
$ws.Cells.Item(7, 16).Value = 0.5796524726490595
$ws.Cells.Item(7, 17).Value = 4.468524906029001
$ws.Cells.Item(7, 18).Value = 40.21672415426101
$ws.Cells.Item(7, 19).Value = 0.2431655199075481
$ws.Cells.Item(7, 20).Value = 0.2431655199075481
$ws.Cells.Item(8, 7).Value = 0.18094
$ws.Cells.Item(8, 8).Value = 0.54282
$ws.Cells.Item(8, 9).Value = 0.1890307689521785
$ws.Cells.Item(8, 10).Value = 0.1890307689521785
$ws.Cells.Item(8, 13).Value = 7.757543333333333
$ws.Cells.Item(8, 14).Value = 23.27263
$ws.Cells.Item(8, 15).Value = 0.4040769763164727
$ws.Cells.Item(8, 16).Value = 0.4040769763164727
$ws.Cells.Item(8, 17).Value = 1.403649890733333
$ws.Cells.Item(8, 18).Value = 12.6328490166
$ws.Cells.Item(8, 19).Value = 0.07638298154897405
$ws.Cells.Item(8, 20).Value = 0.07638298154897405
$ws.Cells.Item(9, 7).Value = 0.18094
$ws.Cells.Item(9, 8).Value = 0.54282
$ws.Cells.Item(9, 9).Value = 0.1890307689521785
$ws.Cells.Item(9, 10).Value = 0.1890307689521785
$ws.Cells.Item(9, 15).Value = 0.01627055103446774
$ws.Cells.Item(9, 16).Value = 0.01627055103446774
$ws.Cells.Item(9, 17).Value = 0.0565193231
$ws.Cells.Item(9, 18).Value = 0.5086739079
$ws.Cells.Item(9, 19).Value = 0.003075634773321101
$ws.Cells.Item(9, 20).Value = 0.003075634773321101
$ws.Cells.Item(10, 7).Value = 0.18094
$ws.Cells.Item(10, 8).Value = 0.54282
$ws.Cells.Item(10, 9).Value = 0.1890307689521785
$ws.Cells.Item(10, 10).Value = 0.1890307689521785
$ws.Cells.Item(10, 13).Value = 11.12827366666667
$ws.Cells.Item(10, 14).Value = 33.384821
$ws.Cells.Item(10, 15).Value = 0.5796524726490594
$ws.Cells.Item(10, 16).Value = 0.5796524726490595
$ws.Cells.Item(10, 17).Value = 2.013549837246666
$ws.Cells.Item(10, 18).Value = 18.12194853522
$ws.Cells.Item(10, 19).Value = 0.1095721526298833
$ws.Cells.Item(10, 20).Value = 0.1095721526298833
